$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3274329
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3274329
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 9822987
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -9823323
$ws.Range("H43").Value = 1267.6
$ws.Range("J43").Value = 1463.5
$ws.Range("L43").Value = 1463.5
$ws.Range("N43").Value = -1601.5
$ws.Range("H129").Value = 855.0928
$ws.Range("J129").Value = 965.4286
$ws.Range("L129").Value = 2896.2858
$ws.Range("N129").Value = -12896.2858
$ws.Range("H132").Value = 1107.2106
$ws.Range("I132").Value = 1107.2106
$ws.Range("K132").Value = 3321.6318
$ws.Range("M132").Value = -791.6318000000001
$ws.Range("H137").Value = 1196.6177
$ws.Range("I137").Value = 1205.7742
$ws.Range("J137").Value = 1102
$ws.Range("K137").Value = 3617.3226
$ws.Range("L137").Value = 3306
$ws.Range("M137").Value = -1067.3226
$ws.Range("N137").Value = -8406
$ws.Range("H138").Value = 2337.2
$ws.Range("I138").Value = 843.91174
$ws.Range("J138").Value = 5510.4375
$ws.Range("K138").Value = 2531.73522
$ws.Range("L138").Value = 16531.3125
$ws.Range("M138").Value = 2608.26478
$ws.Range("N138").Value = -26811.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6459.5
$ws.Range("I32").Value = 4549.1387
$ws.Range("K32").Value = 4549.1387
$ws.Range("M32").Value = -4262.1387
$ws.Range("H45").Value = 8461.076999999999
$ws.Range("I45").Value = 9082.833000000001
$ws.Range("K45").Value = 9082.833000000001
$ws.Range("M45").Value = -8705.833000000001
$ws.Range("H88").Value = 2843.2856
$ws.Range("I88").Value = 2779.2
$ws.Range("J88").Value = 3003.5
$ws.Range("K88").Value = 2779.2
$ws.Range("L88").Value = 3003.5
$ws.Range("M88").Value = -2373.2
$ws.Range("N88").Value = -3815.5
$ws.Range("H91").Value = 2843.2856
$ws.Range("I91").Value = 2779.2
$ws.Range("J91").Value = 3003.5
$ws.Range("K91").Value = 2779.2
$ws.Range("L91").Value = 3003.5
$ws.Range("M91").Value = -1375.2
$ws.Range("N91").Value = -5811.5
$ws.Range("H132").Value = 3148.45
$ws.Range("I132").Value = 1625.2609
$ws.Range("K132").Value = 4875.7827
$ws.Range("M132").Value = -2345.7827
$ws.Range("H134").Value = 41357.25
$ws.Range("J134").Value = 41357.25
$ws.Range("L134").Value = 41357.25
$ws.Range("N134").Value = -51497.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 19426.285
$ws.Range("I20").Value = 1446.8
$ws.Range("J20").Value = 64375
$ws.Range("K20").Value = 1446.8
$ws.Range("L20").Value = 64375
$ws.Range("M20").Value = -1199.8
$ws.Range("N20").Value = -64869
$ws.Range("H61").Value = 20000
$ws.Range("J61").Value = 20000
$ws.Range("L61").Value = 20000
$ws.Range("N61").Value = -20626
$ws.Range("H86").Value = 5954208.5
$ws.Range("J86").Value = 2135.0833
$ws.Range("L86").Value = 2135.0833
$ws.Range("N86").Value = -4381.0833
$ws.Range("H89").Value = 5954208.5
$ws.Range("J89").Value = 2135.0833
$ws.Range("L89").Value = 10675.4165
$ws.Range("N89").Value = -21907.4165
$ws.Range("H99").Value = 55556612
$ws.Range("I99").Value = 76923810
$ws.Range("J99").Value = 1893.8
$ws.Range("K99").Value = 76923810
$ws.Range("L99").Value = 1893.8
$ws.Range("M99").Value = -76922312
$ws.Range("N99").Value = -4889.8
$ws.Range("H105").Value = 2981.8667
$ws.Range("I105").Value = 2108
$ws.Range("K105").Value = 2108
$ws.Range("M105").Value = -361
$ws.Range("H134").Value = 3784.625
$ws.Range("I134").Value = 4060.205
$ws.Range("J134").Value = 2590.4443
$ws.Range("K134").Value = 12180.615
$ws.Range("L134").Value = 7771.3329
$ws.Range("M134").Value = -9645.615
$ws.Range("N134").Value = -12841.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1107.7222
$ws.Range("I16").Value = 662.63635
$ws.Range("J16").Value = 1807.1428
$ws.Range("K16").Value = 662.63635
$ws.Range("L16").Value = 1807.1428
$ws.Range("M16").Value = -375.63635
$ws.Range("N16").Value = -2381.1428
$ws.Range("H31").Value = 3086.606
$ws.Range("I31").Value = 1847.76
$ws.Range("K31").Value = 1847.76
$ws.Range("M31").Value = -1552.76
$ws.Range("H34").Value = 3086.606
$ws.Range("I34").Value = 1847.76
$ws.Range("K34").Value = 1847.76
$ws.Range("M34").Value = -1645.76
$ws.Range("H58").Value = 1640
$ws.Range("I58").Value = 1308.8889
$ws.Range("J58").Value = 2633.3333
$ws.Range("K58").Value = 1308.8889
$ws.Range("L58").Value = 2633.3333
$ws.Range("M58").Value = -1105.8889
$ws.Range("N58").Value = -3039.3333
$ws.Range("H94").Value = 3542.68
$ws.Range("I94").Value = 4000.6
$ws.Range("J94").Value = 3237.4
$ws.Range("K94").Value = 4000.6
$ws.Range("L94").Value = 3237.4
$ws.Range("M94").Value = -3549.6
$ws.Range("N94").Value = -4139.4
$ws.Range("H99").Value = 3196.484
$ws.Range("I99").Value = 2682.8333
$ws.Range("J99").Value = 3907.6924
$ws.Range("K99").Value = 2682.8333
$ws.Range("L99").Value = 3907.6924
$ws.Range("M99").Value = -1184.8333
$ws.Range("N99").Value = -6903.6924
$ws.Range("H107").Value = 242.51282
$ws.Range("I107").Value = 191.18182
$ws.Range("J107").Value = 262.67856
$ws.Range("K107").Value = 191.18182
$ws.Range("L107").Value = 262.67856
$ws.Range("M107").Value = 1728.81818
$ws.Range("N107").Value = -4102.67856
$ws.Range("H113").Value = 1107.7222
$ws.Range("I113").Value = 662.63635
$ws.Range("J113").Value = 1807.1428
$ws.Range("K113").Value = 662.63635
$ws.Range("L113").Value = 1807.1428
$ws.Range("M113").Value = 1507.36365
$ws.Range("N113").Value = -6147.1428
$ws.Range("H126").Value = 3196.484
$ws.Range("I126").Value = 2682.8333
$ws.Range("J126").Value = 3907.6924
$ws.Range("K126").Value = 8048.499899999999
$ws.Range("L126").Value = 11723.0772
$ws.Range("M126").Value = -5578.499899999999
$ws.Range("N126").Value = -16663.0772
$ws.Range("H134").Value = 2964.5
$ws.Range("I134").Value = 3317
$ws.Range("K134").Value = 9951
$ws.Range("M134").Value = -7416
$ws.Range("H136").Value = 1640
$ws.Range("I136").Value = 1308.8889
$ws.Range("J136").Value = 2633.3333
$ws.Range("K136").Value = 3926.6667
$ws.Range("L136").Value = 7899.999899999999
$ws.Range("M136").Value = -1376.6667
$ws.Range("N136").Value = -12999.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 2920.2666
$ws.Range("I94").Value = 624
$ws.Range("J94").Value = 3084.2856
$ws.Range("K94").Value = 1872
$ws.Range("L94").Value = 9252.856800000001
$ws.Range("M94").Value = -1196
$ws.Range("N94").Value = -10604.8568
$ws.Range("H131").Value = 909.25
$ws.Range("I131").Value = 495.27274
$ws.Range("J131").Value = 1047.2424
$ws.Range("K131").Value = 1485.81822
$ws.Range("L131").Value = 3141.7272
$ws.Range("M131").Value = 3554.18178
$ws.Range("N131").Value = -13221.7272

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6277.222
$ws.Range("I70").Value = 6206
$ws.Range("K70").Value = 6206
$ws.Range("M70").Value = -5936
$ws.Range("H73").Value = 6277.222
$ws.Range("I73").Value = 6206
$ws.Range("K73").Value = 6206
$ws.Range("M73").Value = -5270
$ws.Range("H102").Value = 1087.3334
$ws.Range("I102").Value = 841.3333
$ws.Range("K102").Value = 841.3333
$ws.Range("M102").Value = 780.6667
$ws.Range("H126").Value = 5468.385
$ws.Range("I126").Value = 6413.85
$ws.Range("J126").Value = 2316.8333
$ws.Range("K126").Value = 19241.55
$ws.Range("L126").Value = 6950.499899999999
$ws.Range("M126").Value = -16771.55
$ws.Range("N126").Value = -11890.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2045.125
$ws.Range("I61").Value = 1631.9445
$ws.Range("J61").Value = 2576.3572
$ws.Range("K61").Value = 1631.9445
$ws.Range("L61").Value = 2576.3572
$ws.Range("M61").Value = -1429.9445
$ws.Range("N61").Value = -2980.3572
$ws.Range("H113").Value = 2045.125
$ws.Range("I113").Value = 1631.9445
$ws.Range("J113").Value = 2576.3572
$ws.Range("K113").Value = 1631.9445
$ws.Range("L113").Value = 2576.3572
$ws.Range("M113").Value = 538.0554999999999
$ws.Range("N113").Value = -6916.3572
$ws.Range("H134").Value = 43850
$ws.Range("J134").Value = 43850
$ws.Range("L134").Value = 43850
$ws.Range("N134").Value = -53990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 803.5357
$ws.Range("I113").Value = 718.5
$ws.Range("J113").Value = 956.6
$ws.Range("K113").Value = 2155.5
$ws.Range("L113").Value = 2869.8
$ws.Range("M113").Value = 14.5
$ws.Range("N113").Value = -7209.8
$ws.Range("H132").Value = 1299.2979
$ws.Range("I132").Value = 849.1389
$ws.Range("J132").Value = 2772.5454
$ws.Range("K132").Value = 2547.4167
$ws.Range("L132").Value = 8317.636200000001
$ws.Range("M132").Value = -17.41670000000022
$ws.Range("N132").Value = -13377.6362
$ws.Range("H136").Value = 1000.8788
$ws.Range("I136").Value = 593.4231
$ws.Range("J136").Value = 2514.2856
$ws.Range("K136").Value = 1780.2693
$ws.Range("L136").Value = 7542.8568
$ws.Range("M136").Value = 769.7307000000001
$ws.Range("N136").Value = -12642.8568
